$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new row between "Contact" (row 10) and "Description" (row 11),
# shifting "Description" and everything below it down by one row.
$ws.Rows.Item(11).Insert()

# Match the formatting of the other data rows (style used by rows 2-21).
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

# Populate the new "Jurisdiction" property row (value intentionally blank).
$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""

# Refresh the publication date (row 8, column B - "Date" row, unaffected by
# the row insert above it).
$ws.Cells.Item(8, 2).Value = "2024-09-12T14:01:50+00:00"
